# clase4/teorica_4.pptx — "Actualizo ppts 4 y 5"
#
# Slide 6 ("Pensemos los siguientes problemas"), shape "CuadroTexto 4"
# (id=5) is resized/repositioned slightly and four of its bulleted
# paragraphs are switched to justified alignment.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Locate the "CuadroTexto 4" textbox on this slide.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "CuadroTexto 4") {
        $shape = $candidate
    }
}

# Narrow the box (cx 5897527 -> 5432952 EMU) and nudge its left edge by
# a single EMU (127589 -> 127590); top/height are untouched.
$shape.Left = 10.0465
$shape.Width = 427.7915

# Paragraphs 3, 4, 5 and 6 (the marL=177800/indent=-177800 bullet items)
# move to justified alignment (ppAlignJustify = 4).
$tr = $shape.TextFrame.TextRange
$tr.Paragraphs(3, 1).ParagraphFormat.Alignment = 4
$tr.Paragraphs(4, 1).ParagraphFormat.Alignment = 4
$tr.Paragraphs(5, 1).ParagraphFormat.Alignment = 4
$tr.Paragraphs(6, 1).ParagraphFormat.Alignment = 4
